$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values can look numeric (e.g. "1.00", "30.46") and Excel would
# auto-convert them to actual numbers, losing the original text formatting.
# Force them to remain text by setting the NumberFormat to Text ("@") before
# assigning the value, then resetting the style back to Normal so no stray
# style index is left on the cell (matches original formatting exactly).
$priceUpdates = [ordered]@{
    'D2' = '66.448.09'
    'D3' = '3.458.75'
    'D5' = '588.31'
    'D6' = '176.57'
    'D9' = '3.455.37'
    'D13' = '4.056.76'
    'D14' = '30.46'
    'D15' = '0.134'
    'D16' = '66.284.91'
    'D17' = '0.0000173'
    'D18' = '3.453.74'
    'D20' = '13.82'
    'D21' = '374.14'
    'D22' = '7.63'
    'D23' = '73.38'
    'D24' = '0.998'
    'D26' = '0.535'
    'D27' = '9.90'
    'D29' = '1.00'
    'D30' = '5.90'
    'D31' = '2.01'
    'D32' = '23.75'
    'D33' = '1.00'
    'D34' = '7.04'
    'D35' = '1.28'
    'D37' = '160.67'
    'D39' = '28.48'
    'D40' = '1.82'
    'D42' = '4.51'
    'D43' = '2.761.36'
    'D44' = '6.42'
    'D45' = '0.0695'
    'D46' = '25.33'
    'D47' = '339.59'
    'D48' = '40.06'
    'D49' = '0.0293'
    'D51' = '0.993'
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# Coin name / link / percentage-change cells are never numeric-looking
# (letters, URLs, or "%"-suffixed strings), so plain assignment keeps them text.
$textUpdates = [ordered]@{
    'E2' = '  -0.68%  '
    'E3' = '  -1.97%  '
    'E5' = '  -0.03%  '
    'E6' = '  -0.47%  '
    'E7' = '  +1.92%  '
    'E8' = '  -0.02%  '
    'E9' = '  -2.00%  '
    'E10' = '  -1.98%  '
    'E11' = '  +0.22%  '
    'E12' = '  -1.82%  '
    'E13' = '  -1.90%  '
    'B14' = 'Avalanche'
    'C14' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'E14' = '  -0.61%  '
    'B15' = 'TRON'
    'C15' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'E15' = '  +1.22%  '
    'E16' = '  -0.91%  '
    'E17' = '  -1.07%  '
    'E18' = '  -2.01%  '
    'E19' = '  -2.30%  '
    'E20' = '  -1.76%  '
    'E21' = '  -2.35%  '
    'E22' = '  -3.27%  '
    'E23' = '  +2.29%  '
    'E24' = '  -0.40%  '
    'E25' = '  +2.91%  '
    'E26' = '  -0.83%  '
    'E27' = '  -0.58%  '
    'E28' = '  +2.34%  '
    'E29' = '  -0.08%  '
    'E30' = '  -1.41%  '
    'E31' = '  -0.94%  '
    'E32' = '  -3.57%  '
    'E33' = '  +0.01%  '
    'B34' = 'Aptos'
    'C34' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'E34' = '  -3.09%  '
    'B35' = 'Fetch.AI'
    'C35' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'E35' = '  -6.88%  '
    'E36' = '  -1.86%  '
    'E37' = '  +0.86%  '
    'E38' = '  -0.89%  '
    'E39' = '  -4.09%  '
    'E40' = '  +0.60%  '
    'E41' = '  +0.87%  '
    'E42' = '  -0.64%  '
    'E43' = '  +1.27%  '
    'E44' = '  -3.11%  '
    'E45' = '  -1.99%  '
    'E46' = '  -0.79%  '
    'E47' = '  +3.64%  '
    'E48' = '  -1.63%  '
    'E49' = '  -2.33%  '
    'E50' = '  +0.00%  '
    'E51' = '  -3.48%  '
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}
